$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate the sample "class code" data in columns B and C (rows 2-10) ---
# Written in this exact order so the shared-strings table is built up the
# same way Excel would have produced it.
$ws.Range("B2").Value = "sample  B 1"
$ws.Range("B3").Value = "sample  B 2"
$ws.Range("C2").Value = "sample C  1"
$ws.Range("C3").Value = "sample C  2"
$ws.Range("B4").Value = "sample  B 3"
$ws.Range("C4").Value = "sample C  3"
$ws.Range("B5").Value = "sample  B 4"
$ws.Range("C5").Value = "sample C  4"
$ws.Range("B6").Value = "sample  B 5"
$ws.Range("C6").Value = "sample C  5"
$ws.Range("B7").Value = "sample  B 6"
$ws.Range("C7").Value = "sample C  6"
$ws.Range("B8").Value = "sample  B 7"
$ws.Range("C8").Value = "sample C  7"
$ws.Range("B9").Value = "sample  B 8"
$ws.Range("C9").Value = "sample C  8"
$ws.Range("B10").Value = "sample  B 9"
$ws.Range("C10").Value = "sample C  9"

# --- Widen columns B:C to fit the new text, replacing the old bestFit width ---
$ws.Range("B1:C1").ColumnWidth = 22.6

# --- Register the small (size 8) font used for phonetic/info text in the
# workbook's style table, without leaving it applied to any visible cell ---
$ws.Range("Z99").Font.Size = 8
$ws.Range("Z99").Clear()

# --- Move the active selection to match the author's last position ---
$ws.Range("E6").Select() | Out-Null
